$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 11
$ws.Range("H3").Value = 11

$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 2
$ws.Range("H6").Value = 2

$ws.Range("F17").Value = 43
$ws.Range("H17").Value = 43

$ws.Range("E36").Value = 91

$ws.Range("E39").Value = 22

$ws.Range("F42").Value = 13
$ws.Range("H42").Value = 13

$ws.Range("F48").Value = 17
$ws.Range("H48").Value = 17

$ws.Range("E49").Value = 58

$ws.Range("E51").Value = 7

$ws.Range("E54").Value = 1

$ws.Range("E60").Value = 18
$ws.Range("F60").Value = 8
$ws.Range("H60").Value = 8

$ws.Range("F64").Value = 16
$ws.Range("H64").Value = 16
